$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.059880695826386
$ws.Range("D2").Value = 1.066195651312384
$ws.Range("E2").Value = 1.06559725020907
$ws.Range("F2").Value = 1.076294384897703
$ws.Range("J2").Value = 1.064864431662581
$ws.Range("K2").Value = 1.068906852217692
$ws.Range("L2").Value = 1.068310061782959
$ws.Range("M2").Value = 1.078978695818274
$ws.Range("N2").Value = 1.024944988265652

$ws.Range("C3").Value = 1.061845146055403
$ws.Range("D3").Value = 1.068008479248404
$ws.Range("E3").Value = 1.067371612228633
$ws.Range("F3").Value = 1.078204860854018
$ws.Range("J3").Value = 1.066477460189503
$ws.Range("K3").Value = 1.070532738703226
$ws.Range("L3").Value = 1.069897459242695
$ws.Range("M3").Value = 1.080703977836674
$ws.Range("N3").Value = 1.025517505399297

$ws.Range("C4").Value = 1.063113377602309
$ws.Range("D4").Value = 1.069178901813844
$ws.Range("E4").Value = 1.068517146255218
$ws.Range("F4").Value = 1.079438591075603
$ws.Range("J4").Value = 1.067518081017154
$ws.Range("K4").Value = 1.071581758037445
$ws.Range("L4").Value = 1.070921567324905
$ws.Range("M4").Value = 1.081817449544657
$ws.Range("N4").Value = 1.025886005745347

$ws.Range("C5").Value = 1.0636458680607
$ws.Range("D5").Value = 1.06967034282929
$ws.Range("E5").Value = 1.068998124028754
$ws.Range("F5").Value = 1.079956676848267
$ws.Range("J5").Value = 1.067954828209848
$ws.Range("K5").Value = 1.072022054706458
$ws.Range("L5").Value = 1.071351389393251
$ws.Range("M5").Value = 1.082284875486469
$ws.Range("N5").Value = 1.026040460367392

$ws.Range("C6").Value = 1.063735236599869
$ws.Range("D6").Value = 1.069752822985298
$ws.Range("E6").Value = 1.069078847341386
$ws.Range("F6").Value = 1.080043632571592
$ws.Range("J6").Value = 1.068028117659702
$ws.Range("K6").Value = 1.072095941226431
$ws.Range("L6").Value = 1.071423517057892
$ws.Range("M6").Value = 1.082363319116277
$ws.Range("N6").Value = 1.026066366989901

$ws.Range("C7").Value = 1.063120495387431
$ws.Range("D7").Value = 1.069185470824824
$ws.Range("E7").Value = 1.068523575456854
$ws.Range("F7").Value = 1.079445515995303
$ws.Range("J7").Value = 1.067523919699971
$ws.Range("K7").Value = 1.071587644074964
$ws.Range("L7").Value = 1.070927313408353
$ws.Range("M7").Value = 1.08182369795223
$ws.Range("N7").Value = 1.025888071387038

$ws.Range("C8").Value = 1.060545203331271
$ws.Range("D8").Value = 1.066808853530931
$ws.Range("E8").Value = 1.066197452109528
$ws.Range("F8").Value = 1.07694056277061
$ws.Range("J8").Value = 1.065410217344837
$ws.Range("K8").Value = 1.069456967007027
$ws.Range("L8").Value = 1.068847170079769
$ws.Range("M8").Value = 1.079562373911782
$ws.Range("N8").Value = 1.02513888153306

$ws.Range("C9").Value = 1.055984105720611
$ws.Range("D9").Value = 1.062600264868534
$ws.Range("E9").Value = 1.062077872050329
$ws.Range("F9").Value = 1.072506734448363
$ws.Range("J9").Value = 1.061661000500861
$ws.Range("K9").Value = 1.065678442843347
$ws.Range("L9").Value = 1.06515766306946
$ws.Range("M9").Value = 1.075554647605435
$ws.Range("N9").Value = 1.023803481430824

$ws.Range("C10").Value = 1.052926576212342
$ws.Range("D10").Value = 1.059779528477302
$ws.Range("E10").Value = 1.059316524790413
$ws.Range("F10").Value = 1.069536395660844
$ws.Range("J10").Value = 1.059143948768713
$ws.Range("K10").Value = 1.06314226296723
$ws.Range("L10").Value = 1.062680833646137
$ws.Range("M10").Value = 1.072866301648161
$ws.Range("N10").Value = 1.022902625331243

$ws.Range("C11").Value = 1.051598365262403
$ws.Range("D11").Value = 1.058554311315541
$ws.Range("E11").Value = 1.05811704204199
$ws.Range("F11").Value = 1.068246519868483
$ws.Range("J11").Value = 1.058049638298319
$ws.Range("K11").Value = 1.062039767897639
$ws.Range("L11").Value = 1.061604045054462
$ws.Range("M11").Value = 1.071698057242991
$ws.Range("N11").Value = 1.022509954756825

$ws.Range("C12").Value = 1.051104341028799
$ws.Range("D12").Value = 1.05809861587819
$ws.Range("E12").Value = 1.057670908844889
$ws.Range("F12").Value = 1.067766824418724
$ws.Range("J12").Value = 1.057642480039465
$ws.Range("K12").Value = 1.06162958431155
$ws.Range("L12").Value = 1.061203411594795
$ws.Range("M12").Value = 1.071263471498063
$ws.Range("N12").Value = 1.022363702895156

$ws.Range("C13").Value = 1.051210341507673
$ws.Range("D13").Value = 1.058196391377689
$ws.Range("E13").Value = 1.057766633051658
$ws.Range("F13").Value = 1.067869747245153
$ws.Range("J13").Value = 1.057729848082199
$ws.Range("K13").Value = 1.061717600626473
$ws.Range("L13").Value = 1.06128937929498
$ws.Range("M13").Value = 1.071356721265141
$ws.Range("N13").Value = 1.022395092470031

$ws.Range("C14").Value = 1.051557542825474
$ws.Range("D14").Value = 1.05851665565562
$ws.Range("E14").Value = 1.058080176728629
$ws.Range("F14").Value = 1.068206880017271
$ws.Range("J14").Value = 1.058015996511618
$ws.Range("K14").Value = 1.062005875734165
$ws.Range("L14").Value = 1.061570942271986
$ws.Range("M14").Value = 1.071662147585232
$ws.Range("N14").Value = 1.022497873661107

$ws.Range("C15").Value = 1.051771375823083
$ws.Range("D15").Value = 1.058713901529691
$ws.Range("E15").Value = 1.058273282276896
$ws.Range("F15").Value = 1.068414521361466
$ws.Range("J15").Value = 1.058192210945734
$ws.Range("K15").Value = 1.062183402447912
$ws.Range("L15").Value = 1.061744333646822
$ws.Range("M15").Value = 1.071850244359646
$ws.Range("N15").Value = 1.022561147795984

$ws.Range("C16").Value = 1.053014632933413
$ws.Range("D16").Value = 1.059860759811118
$ws.Range("E16").Value = 1.059396048653099
$ws.Range("F16").Value = 1.069621920624345
$ws.Range("J16").Value = 1.059216479925519
$ws.Range("K16").Value = 1.063215339380078
$ws.Range("L16").Value = 1.062752204201316
$ws.Range("M16").Value = 1.072943744492237
$ws.Range("N16").Value = 1.022928630352987

$ws.Range("C17").Value = 1.053793332094963
$ws.Range("D17").Value = 1.060579116072282
$ws.Range("E17").Value = 1.060099298045009
$ws.Range("F17").Value = 1.070378284137798
$ws.Range("J17").Value = 1.059857781946441
$ws.Range("K17").Value = 1.063861477596405
$ws.Range("L17").Value = 1.06338324852623
$ws.Range("M17").Value = 1.073628536416645
$ws.Range("N17").Value = 1.023158443197648

$ws.Range("C18").Value = 1.054247122801946
$ws.Range("D18").Value = 1.060997754037345
$ws.Range("E18").Value = 1.060509126179034
$ws.Range("F18").Value = 1.070819102753552
$ws.Range("J18").Value = 1.060231418383327
$ws.Range("K18").Value = 1.064237944356233
$ws.Range("L18").Value = 1.063750911949205
$ws.Range("M18").Value = 1.074027562759262
$ws.Range("N18").Value = 1.023292239329577

$ws.Range("C19").Value = 1.05440178467115
$ws.Range("D19").Value = 1.061140437142803
$ws.Range("E19").Value = 1.060648805640094
$ws.Range("F19").Value = 1.070969350895387
$ws.Range("J19").Value = 1.060358747405625
$ws.Range("K19").Value = 1.064366240046564
$ws.Range("L19").Value = 1.063876206028758
$ws.Range("M19").Value = 1.074163553006094
$ws.Range("N19").Value = 1.023337818185288

$ws.Range("C20").Value = 1.053709827767501
$ws.Range("D20").Value = 1.060502081367641
$ws.Range("E20").Value = 1.060023883979475
$ws.Range("F20").Value = 1.070297170415312
$ws.Range("J20").Value = 1.059789020328171
$ws.Range("K20").Value = 1.063792196121959
$ws.Range("L20").Value = 1.06331558641015
$ws.Range("M20").Value = 1.073555106372267
$ws.Range("N20").Value = 1.023133812342748

$ws.Range("C21").Value = 1.051455319316628
$ws.Range("D21").Value = 1.058422362412337
$ws.Range("E21").Value = 1.057987862438358
$ws.Range("F21").Value = 1.068107618959236
$ws.Range("J21").Value = 1.057931751936929
$ws.Range("K21").Value = 1.061921004496671
$ws.Range("L21").Value = 1.061488047541638
$ws.Range("M21").Value = 1.071572225196358
$ws.Range("N21").Value = 1.022467618163367

$ws.Range("C22").Value = 1.050033945810124
$ws.Range("D22").Value = 1.057111305962142
$ws.Range("E22").Value = 1.056704299849442
$ws.Range("F22").Value = 1.066727604795402
$ws.Range("J22").Value = 1.05676005402948
$ws.Range("K22").Value = 1.060740638148381
$ws.Range("L22").Value = 1.060335136896346
$ws.Range("M22").Value = 1.070321750871318
$ws.Range("N22").Value = 1.022046457816724

$ws.Range("C23").Value = 1.050787818305922
$ws.Range("D23").Value = 1.057806656451158
$ws.Range("E23").Value = 1.05738507328281
$ws.Range("F23").Value = 1.067459502080612
$ws.Range("J23").Value = 1.057381575323464
$ws.Range("K23").Value = 1.061366746604084
$ws.Range("L23").Value = 1.060946689499404
$ws.Range("M23").Value = 1.070985014230036
$ws.Range("N23").Value = 1.022269942973165

$ws.Range("C24").Value = 1.053747561049085
$ws.Range("D24").Value = 1.060536891168556
$ws.Range("E24").Value = 1.060057961477949
$ws.Range("F24").Value = 1.070333823312749
$ws.Range("J24").Value = 1.059820092053893
$ws.Range("K24").Value = 1.063823502719113
$ws.Range("L24").Value = 1.063346161286894
$ws.Range("M24").Value = 1.07358828748778
$ws.Range("N24").Value = 1.023144942737654

$ws.Range("C25").Value = 1.057166125959339
$ws.Range("D25").Value = 1.063690850005484
$ws.Range("E25").Value = 1.063145439626397
$ws.Range("F25").Value = 1.073655447433862
$ws.Range("J25").Value = 1.062633284044242
$ws.Range("K25").Value = 1.066658231033261
$ws.Range("L25").Value = 1.066114439089634
$ws.Range("M25").Value = 1.076593574600591
$ws.Range("N25").Value = 1.024150555004769
